$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow {
    param($row, $name, $b, $c, $d, $e, $f, $g, $h)
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- El Salvador's stats were updated, which pushed it above Maldivas
#     in the list (sorted descending by total cases); the two rows swap
#     country names while Maldivas keeps its previous figures. ---
Set-CountryRow 114 "El Salvador" 555 65 180 363 3 1 12
Set-CountryRow 115 "Maldivas" 527 0 18 508 2 0 1

# --- Timor Oriental (row 183): redistribute recuperados/activos ---
$ws.Cells.Item(183, 4).Value = 20
$ws.Cells.Item(183, 5).Value = 4

# --- Belice / Santa Lucia swap (rows 188, 189) ---
Set-CountryRow 188 "Belice" 18 0 13 3 1 0 2
Set-CountryRow 189 "Santa Lucia" 18 0 15 3 0 0 0

# --- San Vicente y las Granadinas / Namibia swap (rows 194, 195) ---
Set-CountryRow 194 "San Vicente y las Granadinas" 16 0 8 8 0 0 0
Set-CountryRow 195 "Namibia" 16 0 8 8 0 0 0

# --- San Cristobal y Nieves / Burundi swap (rows 198, 199) ---
Set-CountryRow 198 "San Cristobal y Nieves" 15 0 8 7 0 0 0
Set-CountryRow 199 "Burundi" 15 0 7 7 0 0 1

# --- Papua Nueva Guinea (row 209): redistribute recuperados/activos ---
$ws.Cells.Item(209, 4).Value = 6
$ws.Cells.Item(209, 5).Value = 2
